$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Rename the Pearson logo pictures embedded in both footers
# (PNG, currently "image1.png") to "image2.png".
for ($i = 1; $i -le 2; $i++) {
    $footer = $sec.Footers($i)
    if ($footer.Exists) {
        $pic = $footer.Range.InlineShapes(1)
        $shape = $pic.ConvertToShape()
        $shape.Name = "image2.png"
        [void]$shape.ConvertToInlineShape()
    }
}

# Rename the BTEC logo pictures embedded in both headers
# (JPEG, currently "image2.jpg") to "image1.jpg".
for ($i = 1; $i -le 2; $i++) {
    $header = $sec.Headers($i)
    if ($header.Exists) {
        $pic = $header.Range.InlineShapes(1)
        $shape = $pic.ConvertToShape()
        $shape.Name = "image1.jpg"
        [void]$shape.ConvertToInlineShape()
    }
}
